$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Bump "Last Updated" (column E) from 2026-01-18 -> 2026-01-20 for every
#    data row that currently carries that date. A handful of rows (25, 32,
#    47, 67, 79, 84) were already stamped 2026-01-17 and are left untouched.
#    We route the write through Formula + Copy/PasteSpecial(values) so the
#    result lands back as a plain text value (matching the original
#    inlineStr/text cell) instead of Excel's automatic "looks like a date"
#    number coercion, while keeping the existing cell style untouched.
# ---------------------------------------------------------------------------
$skipRows = @(25, 32, 47, 67, 79, 84)

for ($r = 2; $r -le 88; $r++) {
    if ($skipRows -contains $r) { continue }
    $cell = $ws.Cells.Item($r, 5)   # column E
    $cell.Formula = '="2026-01-20"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)       # xlPasteValues
}
$excel.CutCopyMode = $false         # drop the clipboard so it can't bleed into later inserts

# ---------------------------------------------------------------------------
# 2) Row 89 used to be fully blank; it now holds a new competitor record.
# ---------------------------------------------------------------------------
$ws.Range("A89").Value = "Phreesia Audit"
$ws.Range("B89").Value = "https://phreesia.com"
$ws.Range("C89").Value = "Active"
$ws.Range("D89").Value = "Medium"

$e89 = $ws.Cells.Item(89, 5)
$e89.Formula = '="2026-01-20"'
$e89.Copy()
$e89.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) The sheet gained one more trailing blank row (98), pushing the used
#    range/dimension from A1:AF97 to A1:AF98. Inserting a row right after
#    the current last row (97) appends a new, fully empty styled row
#    without disturbing any existing data above it.
# ---------------------------------------------------------------------------
$ws.Rows.Item(98).Insert()
